$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for "Primera" quality on 2022-07-08
# (serial 44750). This pushes every subsequent record down by one row
# (old row 5 -> new row 6, ..., old row 38 -> new row 39).
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = 44750
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 100112043
$ws.Range("G5").Value = "Pepino dulce"
$ws.Range("H5").Value = "Cultivar IV Región"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("N5").Value = "`$/bandeja 18 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 833
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"
